# Apply cryptos list update (price & volume refresh + one pair of rows reordered)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.478.68"
$ws.Range("E2").Value = "  +4.28%  "
$ws.Range("D3").Value = "2.426.46"
$ws.Range("E3").Value = "  +5.16%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.54"
$ws.Range("E5").Value = "  +2.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.85"
$ws.Range("E6").Value = "  +7.31%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  +2.22%  "
$ws.Range("D9").Value = "2.425.09"
$ws.Range("E9").Value = "  +5.16%  "
$ws.Range("E10").Value = "  +3.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.77"
$ws.Range("E11").Value = "  +3.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.151"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  +5.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.17"
$ws.Range("E14").Value = "  +12.10%  "
$ws.Range("D15").Value = "2.858.96"
$ws.Range("E15").Value = "  +5.17%  "
$ws.Range("D16").Value = "62.332.16"
$ws.Range("E16").Value = "  +4.10%  "
$ws.Range("E17").Value = "  +7.40%  "
$ws.Range("D18").Value = "2.421.99"
$ws.Range("E18").Value = "  +5.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.17"
$ws.Range("E19").Value = "  +6.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.88"
$ws.Range("E20").Value = "  +10.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.21"
$ws.Range("E21").Value = "  +3.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.83"
$ws.Range("E22").Value = "  +4.37%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.18"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.55"
$ws.Range("E27").Value = "  +14.98%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("E29").Value = "  +5.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.36"
$ws.Range("E30").Value = "  +16.35%  "
$ws.Range("E31").Value = "  +5.33%  "
$ws.Range("D32").Value = "0.0₃0786"
$ws.Range("E32").Value = "  +8.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.49"
$ws.Range("E33").Value = "  +11.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "172.02"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.43"
$ws.Range("E35").Value = "  +5.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.396"
$ws.Range("E36").Value = "  +4.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "378.47"
$ws.Range("E37").Value = "  +19.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.56"
$ws.Range("E38").Value = "  +5.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.45"
$ws.Range("E39").Value = "  +11.67%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +12.13%  "
$ws.Range("E43").Value = "  +3.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "144.95"
$ws.Range("E44").Value = "  +6.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.66"
$ws.Range("E45").Value = "  +7.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "20.76"
$ws.Range("E46").Value = "  +10.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.590"
$ws.Range("E47").Value = "  +4.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0953"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0521"
$ws.Range("E49").Value = "  +6.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0223"
$ws.Range("E50").Value = "  +4.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.85"
$ws.Range("E51").Value = "  +6.37%  "
